$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the header row cells to match the updated CSV output column names
$ws.Range("A1").Value = "Service Name"
$ws.Range("F1").Value = "Phone Number"
$ws.Range("C1").Value = "Membership"

# Widen column F slightly so the new "Phone Number" header fits
$ws.Columns.Item(6).ColumnWidth = 16.3

# Update the active selection to match the new state (cell C1 selected)
$ws.Range("C1").Select()
